# Updates cryptocurrency price/volume data per the latest GitHub Actions scrape.
# (rows 39/40 also swap Coin/Link since VeChain and FraxShare traded ranks.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.112.73"
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.665.44"
$ws.Range("E3").Value = "  -0.69%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.27"
$ws.Range("E5").Value = "  -0.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5208"
$ws.Range("E6").Value = "  -1.25%  "

$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2596"
$ws.Range("E8").Value = "  -2.90%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06325"
$ws.Range("E9").Value = "  +0.78%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.97"
$ws.Range("E10").Value = "  -1.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07536"
$ws.Range("E11").Value = "  +0.15%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.673.13"
$ws.Range("E12").Value = "  -0.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.409"
$ws.Range("E13").Value = "  -1.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5370"
$ws.Range("E14").Value = "  -4.66%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7966"
$ws.Range("E15").Value = "  -1.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.03"
$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.150.02"
$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("E18").Value = "  -0.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.707"
$ws.Range("E19").Value = "  -2.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "187.40"
$ws.Range("E20").Value = "  -0.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.22"
$ws.Range("E21").Value = "  -2.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.201"
$ws.Range("E22").Value = "  +0.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("E23").Value = "  -0.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.02"
$ws.Range("E24").Value = "  +0.72%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1225"
$ws.Range("E25").Value = "  -2.63%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.408"
$ws.Range("E26").Value = "  -2.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.68"
$ws.Range("E27").Value = "  -0.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06181"
$ws.Range("E28").Value = "  -4.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.361"
$ws.Range("E29").Value = "  +1.91%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.266"
$ws.Range("E30").Value = "  -0.98%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.474"

$ws.Range("E32").Value = "  -2.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.635"
$ws.Range("E33").Value = "  -0.42%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9898"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.391"
$ws.Range("E35").Value = "  -1.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.756"
$ws.Range("E36").Value = "  +1.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5892"
$ws.Range("E37").Value = "  -2.96%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.107.46"
$ws.Range("E38").Value = "  +0.73%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01596"
$ws.Range("E39").Value = "  -0.49%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.019"
$ws.Range("E40").Value = "  -1.96%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8495"
$ws.Range("E41").Value = "  -1.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.003"
$ws.Range("E42").Value = "  -0.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.91"
$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.817.43"
$ws.Range("E44").Value = "  -0.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₈110"
$ws.Range("E45").Value = "  +1.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.33"
$ws.Range("E46").Value = "  -2.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.005"
$ws.Range("E47").Value = "  +0.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.057"
$ws.Range("E48").Value = "  +1.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05247"
$ws.Range("E49").Value = "  -0.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4258"
$ws.Range("E50").Value = "  -0.32%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.873"
$ws.Range("E51").Value = "  -0.98%  "
